$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    3885.644404753759,
    3885.644404753759,
    3885.644404753759,
    3885.644404753759,
    3885.644404753759,
    3742.104510259517,
    3742.104510259517,
    3742.104510259517,
    3742.104510259517,
    3742.104510259517,
    3730.524520977293
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
